$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the "55B / Incluir sitemap e analytics / 30" row from the
#     "PARA FAZER:" (to-do) list (row 99) up into the "FEITO:" (done) list,
#     right after the last done item (row 93). This is modeled as:
#       a) delete row 99 (the 55B row), which shifts rows 100-103 up to 99-102
#       b) insert a new blank row at 94, which shifts everything from 94
#          down by one, restoring the original row numbers for the
#          "FAZENDO:" header / item and "PARA FAZER:" header / remaining
#          items, and recreating a fresh blank row 94 for the moved item.

$ws.Rows("99").Delete() | Out-Null
$ws.Rows("94").Insert() | Out-Null

# --- 2. Populate the newly-inserted row 94 with the moved task, matching
#     the centered-alignment look of the other "col A" entries.
$ws.Range("A94").Value = "55B"
$ws.Range("B94").Value = " Incluir sitemap e analytics "
$ws.Range("C94").Value = 30
$ws.Range("A94").HorizontalAlignment = $ws.Range("A93").HorizontalAlignment

# --- 3. Update the (now single) "FAZENDO:" item, which after the shift
#     above lives at row 97. Re-word it and update the time spent.
$ws.Range("B97").Value = " Tirar gambiarra de exibição de mês INCLUIR descritivo"
$ws.Range("C97").Value = 1920

# --- 4. Update the window's scroll position / selected cell to match
#     where the author ended up after editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 79
$win.ScrollColumn = 1
$ws.Range("F94").Select() | Out-Null
